$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value would otherwise be auto-detected as a
# number by the Value setter (e.g. "306.43") are first switched to
# Text format so the literal string is preserved, matching the source
# data which stores every Price/Volume value as text.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.227.89"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.901.68"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "306.43"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.5351"
$ws.Range("E7").Value = "  +2.48%  "
$ws.Range("D8").Value = "0.3815"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D10").Value = "22.03"
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("D11").Value = "0.9007"
$ws.Range("D12").Value = "0.08190"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "95.74"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").Value = "5.336"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "1.003"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "14.81"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("D17").Value = "0.000008638"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "27.262.67"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "5.026"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("B21").Value = "WrappedEther"
$ws.Range("C21").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D21").Value = "1.114.76"
$ws.Range("E21").Value = "  -41.31%  "
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("D24").Value = "149.62"
$ws.Range("E24").Value = "  +1.42%  "
$ws.Range("D25").Value = "2.287"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").Value = "18.35"
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "1.739"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("D28").Value = "116.64"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").Value = "4.807"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").Value = "4.783"
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("D31").Value = "0.09250"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "0.8294"
$ws.Range("E32").Value = "  +3.52%  "
$ws.Range("D33").Value = "0.05052"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "1.220"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").Value = "3.000"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("D36").Value = "3.328"
$ws.Range("E36").Value = "  -3.37%  "
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("D38").Value = "0.5729"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "0.02001"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").Value = "9.304"
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("D42").Value = "6.574"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "117.13"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "0.4934"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "10.08"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "1.634"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "38.31"
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("D50").Value = "0.06136"
$ws.Range("E50").Value = "  +3.17%  "
$ws.Range("D51").Value = "62.86"
$ws.Range("E51").Value = "  -1.71%  "
